$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 941.0833
$ws.Range("I12").Value = 529.3
$ws.Range("K12").Value = 529.3
$ws.Range("M12").Value = -359.3

$ws.Range("H40").Value = 7310.375
$ws.Range("J40").Value = 9994.666999999999
$ws.Range("L40").Value = 9994.666999999999
$ws.Range("N40").Value = -10344.667

$ws.Range("H74").Value = 16923.076
$ws.Range("I74").Value = 17500
$ws.Range("K74").Value = 17500
$ws.Range("M74").Value = -16564

$ws.Range("H76").Value = 3744.2222
$ws.Range("I76").Value = 3814
$ws.Range("J76").Value = 3500
$ws.Range("K76").Value = 3814
$ws.Range("L76").Value = 3500
$ws.Range("M76").Value = -3499
$ws.Range("N76").Value = -4130

$ws.Range("H77").Value = 16923.076
$ws.Range("I77").Value = 17500
$ws.Range("K77").Value = 87500
$ws.Range("M77").Value = -82820

$ws.Range("H79").Value = 3744.2222
$ws.Range("I79").Value = 3814
$ws.Range("J79").Value = 3500
$ws.Range("K79").Value = 3814
$ws.Range("L79").Value = 3500
$ws.Range("M79").Value = -2722
$ws.Range("N79").Value = -5684

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("N121").Value = 0

$ws.Range("H132").Value = 5396.241
$ws.Range("I132").Value = 5640.407
$ws.Range("J132").Value = 2100
$ws.Range("K132").Value = 16921.221
$ws.Range("L132").Value = 6300
$ws.Range("M132").Value = -14391.221
$ws.Range("N132").Value = -11360

$ws.Range("H137").Value = 2944716
$ws.Range("I137").Value = 3848047
$ws.Range("J137").Value = 8891
$ws.Range("K137").Value = 11544141
$ws.Range("L137").Value = 26673
$ws.Range("M137").Value = -11541591
$ws.Range("N137").Value = -31773

$ws.Range("H138").Value = 2961.016
$ws.Range("I138").Value = 2608.1304
$ws.Range("J138").Value = 3169.1282
$ws.Range("K138").Value = 7824.3912
$ws.Range("L138").Value = 9507.384600000001
$ws.Range("M138").Value = -2684.3912
$ws.Range("N138").Value = -19787.3846

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 226.5
$ws.Range("I5").Value = 90
$ws.Range("J5").Value = 499.5
$ws.Range("K5").Value = 90
$ws.Range("L5").Value = 499.5
$ws.Range("M5").Value = 22
$ws.Range("N5").Value = -723.5

$ws.Range("H32").Value = 2733596.8
$ws.Range("I32").Value = 1279642.5
$ws.Range("J32").Value = 23815936
$ws.Range("K32").Value = 1279642.5
$ws.Range("L32").Value = 23815936
$ws.Range("M32").Value = -1279355.5
$ws.Range("N32").Value = -23816510

$ws.Range("H61").Value = 2228.2415
$ws.Range("I61").Value = 1701.1666
$ws.Range("J61").Value = 3090.7273
$ws.Range("K61").Value = 1701.1666
$ws.Range("L61").Value = 3090.7273
$ws.Range("M61").Value = -1489.1666
$ws.Range("N61").Value = -3514.7273

$ws.Range("H97").Value = 1291.3
$ws.Range("I97").Value = 1291.3
$ws.Range("K97").Value = 1291.3
$ws.Range("M97").Value = -795.3

$ws.Range("H132").Value = 2749.4482
$ws.Range("I132").Value = 2486.75
$ws.Range("K132").Value = 7460.25
$ws.Range("M132").Value = -4930.25

$ws.Range("H136").Value = 2228.2415
$ws.Range("I136").Value = 1701.1666
$ws.Range("J136").Value = 3090.7273
$ws.Range("K136").Value = 5103.4998
$ws.Range("L136").Value = 9272.1819
$ws.Range("M136").Value = -2553.4998
$ws.Range("N136").Value = -14372.1819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 226.5
$ws.Range("I4").Value = 90
$ws.Range("J4").Value = 499.5
$ws.Range("K4").Value = 90
$ws.Range("L4").Value = 499.5
$ws.Range("M4").Value = 25
$ws.Range("N4").Value = -729.5

$ws.Range("H94").Value = 133334310
$ws.Range("I94").Value = 153846980
$ws.Range("K94").Value = 153846980
$ws.Range("M94").Value = -153846529

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2019702.4
$ws.Range("I31").Value = 3005
$ws.Range("K31").Value = 3005
$ws.Range("M31").Value = -2710

$ws.Range("H34").Value = 2019702.4
$ws.Range("I34").Value = 3005
$ws.Range("K34").Value = 3005
$ws.Range("M34").Value = -2803

$ws.Range("H58").Value = 2421.3044
$ws.Range("I58").Value = 1413.6
$ws.Range("J58").Value = 4310.75
$ws.Range("K58").Value = 1413.6
$ws.Range("L58").Value = 4310.75
$ws.Range("M58").Value = -1210.6
$ws.Range("N58").Value = -4716.75

$ws.Range("H86").Value = 7003.294
$ws.Range("I86").Value = 7036.8
$ws.Range("J86").Value = 6752
$ws.Range("K86").Value = 7036.8
$ws.Range("L86").Value = 6752
$ws.Range("M86").Value = -5913.8
$ws.Range("N86").Value = -8998

$ws.Range("H89").Value = 7003.294
$ws.Range("I89").Value = 7036.8
$ws.Range("J89").Value = 6752
$ws.Range("K89").Value = 35184
$ws.Range("L89").Value = 33760
$ws.Range("M89").Value = -29568
$ws.Range("N89").Value = -44992

$ws.Range("H107").Value = 2500922.5
$ws.Range("I107").Value = 4167090.2
$ws.Range("K107").Value = 4167090.2
$ws.Range("M107").Value = -4165170.2

$ws.Range("H122").Value = 1547.625
$ws.Range("I122").Value = 1557.5333
$ws.Range("K122").Value = 4672.5999
$ws.Range("M122").Value = -2222.5999

$ws.Range("H132").Value = 14497670
$ws.Range("I132").Value = 3342.2666
$ws.Range("J132").Value = 41674536
$ws.Range("K132").Value = 10026.7998
$ws.Range("L132").Value = 125023608
$ws.Range("M132").Value = -7496.799800000001
$ws.Range("N132").Value = -125028668

$ws.Range("H134").Value = 2545.05
$ws.Range("I134").Value = 2320.7354
$ws.Range("K134").Value = 6962.206200000001
$ws.Range("M134").Value = -4427.206200000001

$ws.Range("H136").Value = 2421.3044
$ws.Range("I136").Value = 1413.6
$ws.Range("J136").Value = 4310.75
$ws.Range("K136").Value = 4240.799999999999
$ws.Range("L136").Value = 12932.25
$ws.Range("M136").Value = -1690.799999999999
$ws.Range("N136").Value = -18032.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 559.56525
$ws.Range("J12").Value = 602.2778
$ws.Range("L12").Value = 1806.8334
$ws.Range("N12").Value = -2152.8334

$ws.Range("H86").Value = 2057.353
$ws.Range("J86").Value = 2297.3333
$ws.Range("L86").Value = 6891.999899999999
$ws.Range("N86").Value = -9263.999899999999

$ws.Range("H89").Value = 2057.353
$ws.Range("J89").Value = 2297.3333
$ws.Range("L89").Value = 20675.9997
$ws.Range("N89").Value = -32531.9997

$ws.Range("H131").Value = 6587111
$ws.Range("I131").Value = 17875748
$ws.Range("J131").Value = 2072.5
$ws.Range("K131").Value = 53627244
$ws.Range("L131").Value = 6217.5
$ws.Range("M131").Value = -53622204
$ws.Range("N131").Value = -16297.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 316.66666
$ws.Range("I2").Value = 316.66666
$ws.Range("K2").Value = 316.66666
$ws.Range("M2").Value = -203.66666

$ws.Range("H21").Value = 24001.2
$ws.Range("J21").Value = 24001.2
$ws.Range("L21").Value = 24001.2
$ws.Range("N21").Value = -24347.2

$ws.Range("H30").Value = 24001.2
$ws.Range("J30").Value = 24001.2
$ws.Range("L30").Value = 24001.2
$ws.Range("N30").Value = -24211.2

$ws.Range("H132").Value = 2575.0625
$ws.Range("I132").Value = 2456.4285
$ws.Range("K132").Value = 7369.2855
$ws.Range("M132").Value = -4839.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 2201.875
$ws.Range("I9").Value = 212
$ws.Range("J9").Value = 5518.3335
$ws.Range("K9").Value = 212
$ws.Range("L9").Value = 5518.3335
$ws.Range("M9").Value = 12
$ws.Range("N9").Value = -5966.3335

$ws.Range("H22").Value = 5326.25
$ws.Range("I22").Value = 2101.6667
$ws.Range("K22").Value = 2101.6667
$ws.Range("M22").Value = -1806.6667

$ws.Range("H27").Value = 5326.25
$ws.Range("I27").Value = 2101.6667
$ws.Range("K27").Value = 2101.6667
$ws.Range("M27").Value = -1994.6667

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").ClearContents()
$ws.Range("N87").Value = 0

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").ClearContents()
$ws.Range("N90").Value = 0

$ws.Range("H132").Value = 4778.7715
$ws.Range("I132").Value = 1771.12
$ws.Range("J132").Value = 12297.9
$ws.Range("K132").Value = 5313.36
$ws.Range("L132").Value = 36893.7
$ws.Range("M132").Value = -2783.36
$ws.Range("N132").Value = -41953.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 40395.4
$ws.Range("I43").Value = 22027
$ws.Range("J43").Value = 44987.5
$ws.Range("K43").Value = 22027
$ws.Range("L43").Value = 44987.5
$ws.Range("M43").Value = -21878
$ws.Range("N43").Value = -45285.5

$ws.Range("H100").Value = 71430120
$ws.Range("I100").Value = 1847.8889
$ws.Range("K100").Value = 3695.7778
$ws.Range("M100").Value = -3154.7778

$ws.Range("H132").Value = 2435.8914
$ws.Range("I132").Value = 2222.7856
$ws.Range("J132").Value = 4673.5
$ws.Range("K132").Value = 6668.3568
$ws.Range("L132").Value = 14020.5
$ws.Range("M132").Value = -4138.3568
$ws.Range("N132").Value = -19080.5
